$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "26.406.99"
$ws.Range("E2").Value2 = "  -1.00%  "
$ws.Range("D3").Value2 = "1.840.19"
$ws.Range("E3").Value2 = "  -1.35%  "
$ws.Range("D4").Value2 = "'1.001"
$ws.Range("E4").Value2 = "  -0.15%  "
$ws.Range("D5").Value2 = "'260.56"
$ws.Range("E5").Value2 = "  -5.53%  "
$ws.Range("D6").Value2 = "'1.001"
$ws.Range("E6").Value2 = "  -0.08%  "
$ws.Range("D7").Value2 = "'0.5178"
$ws.Range("E7").Value2 = "  -2.03%  "
$ws.Range("D8").Value2 = "'0.3268"
$ws.Range("E8").Value2 = "  -4.20%  "
$ws.Range("D9").Value2 = "'0.06765"
$ws.Range("E9").Value2 = "  -1.78%  "
$ws.Range("D10").Value2 = "'18.66"
$ws.Range("E10").Value2 = "  -6.40%  "
$ws.Range("D11").Value2 = "'0.7758"
$ws.Range("E11").Value2 = "  -2.76%  "
$ws.Range("D12").Value2 = "'0.07753"
$ws.Range("E12").Value2 = "  +0.30%  "
$ws.Range("D13").Value2 = "1.839.47"
$ws.Range("E13").Value2 = "  -2.00%  "
$ws.Range("D14").Value2 = "'87.66"
$ws.Range("E14").Value2 = "  -2.41%  "
$ws.Range("D15").Value2 = "'4.992"
$ws.Range("E15").Value2 = "  -2.95%  "
$ws.Range("D16").Value2 = "'1.001"
$ws.Range("E16").Value2 = "  -0.17%  "
$ws.Range("D17").Value2 = "'13.90"
$ws.Range("E17").Value2 = "  -4.09%  "
$ws.Range("D18").Value2 = "'1.000"
$ws.Range("E18").Value2 = "  -0.13%  "
$ws.Range("D19").Value2 = "'0.000007930"
$ws.Range("E19").Value2 = "  -0.76%  "
$ws.Range("D20").Value2 = "26.454.14"
$ws.Range("E20").Value2 = "  -1.00%  "
$ws.Range("D21").Value2 = "2.077.04"
$ws.Range("E21").Value2 = "  -1.65%  "
$ws.Range("D22").Value2 = "'4.614"
$ws.Range("E22").Value2 = "  -2.58%  "
$ws.Range("D23").Value2 = "'9.522"
$ws.Range("E23").Value2 = "  -4.68%  "
$ws.Range("D24").Value2 = "'5.994"
$ws.Range("E24").Value2 = "  -2.59%  "
$ws.Range("D25").Value2 = "'145.57"
$ws.Range("E25").Value2 = "  -0.22%  "
$ws.Range("D26").Value2 = "'2.188"
$ws.Range("E26").Value2 = "  -7.03%  "
$ws.Range("D27").Value2 = "'1.651"
$ws.Range("E27").Value2 = "  +0.07%  "
$ws.Range("D28").Value2 = "'16.97"
$ws.Range("E28").Value2 = "  -1.79%  "
$ws.Range("D29").Value2 = "'111.37"
$ws.Range("E29").Value2 = "  -1.28%  "
$ws.Range("D30").Value2 = "'4.187"
$ws.Range("E30").Value2 = "  -2.92%  "
$ws.Range("D31").Value2 = "'4.119"
$ws.Range("E31").Value2 = "  -4.70%  "
$ws.Range("D32").Value2 = "'0.08703"
$ws.Range("E32").Value2 = "  -1.90%  "
$ws.Range("D33").Value2 = "'0.04816"
$ws.Range("E33").Value2 = "  -2.24%  "
$ws.Range("D34").Value2 = "'1.128"
$ws.Range("E34").Value2 = "  -2.60%  "
$ws.Range("D35").Value2 = "'0.7182"
$ws.Range("E35").Value2 = "  -0.97%  "
$ws.Range("D36").Value2 = "'2.848"
$ws.Range("E36").Value2 = "  -1.41%  "
$ws.Range("D37").Value2 = "'3.087"
$ws.Range("E37").Value2 = "  -4.81%  "
$ws.Range("D38").Value2 = "'0.01777"
$ws.Range("E38").Value2 = "  -3.98%  "
$ws.Range("D39").Value2 = "'2.216"
$ws.Range("E39").Value2 = "  -4.16%  "
$ws.Range("D40").Value2 = "'0.4822"
$ws.Range("E40").Value2 = "  -5.69%  "
$ws.Range("D41").Value2 = "'111.81"
$ws.Range("E41").Value2 = "  -3.74%  "
$ws.Range("D42").Value2 = "'0.8996"
$ws.Range("E42").Value2 = "  -4.44%  "
$ws.Range("D43").Value2 = "'6.082"
$ws.Range("E43").Value2 = "  -0.82%  "
$ws.Range("E44").Value2 = "  -0.06%  "
$ws.Range("D45").Value2 = "'7.711"
$ws.Range("E45").Value2 = "  -4.18%  "
$ws.Range("D46").Value2 = "'0.05962"
$ws.Range("E46").Value2 = "  -0.64%  "
$ws.Range("D47").Value2 = "'0.4141"
$ws.Range("E47").Value2 = "  -6.25%  "
$ws.Range("D48").Value2 = "'9.016"
$ws.Range("E48").Value2 = "  -2.94%  "
$ws.Range("D49").Value2 = "'35.00"
$ws.Range("E49").Value2 = "  -3.36%  "
$ws.Range("D50").Value2 = "'0.1215"
$ws.Range("E50").Value2 = "  -8.84%  "
$ws.Range("D51").Value2 = "'0.8853"
$ws.Range("E51").Value2 = "  +0.42%  "
